$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates (rich-text shared strings; run formatting is uniform so plain text set is equivalent)
$ws.Range("A8").Value = "Volume 32   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/24/2025  Through  11/30/2025"

# Crime-stat table updates (rows 14-30)
$ws.Range("N14").Value = -90.909090909090

$ws.Range("C15").Value = 4
$ws.Range("F15").Value = 4
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 26
$ws.Range("K15").Value = 100
$ws.Range("L15").Value = 225
$ws.Range("M15").Value = 271.428571428571
$ws.Range("N15").Value = 23.809523809523

$ws.Range("D16").Value = 1
$ws.Range("J16").Value = 128
$ws.Range("K16").Value = -23.4375
$ws.Range("L16").Value = -31.468531468531
$ws.Range("M16").Value = -36.363636363636
$ws.Range("N16").Value = -82.897033158813

$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = -60.869565217391
$ws.Range("I17").Value = 201
$ws.Range("J17").Value = 235
$ws.Range("K17").Value = -14.468085106383
$ws.Range("L17").Value = 6.349206349206
$ws.Range("M17").Value = 79.464285714285
$ws.Range("N17").Value = -19.6

$ws.Range("D18").Value = 1
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 6
$ws.Range("I18").Value = 116
$ws.Range("J18").Value = 148
$ws.Range("K18").Value = -21.621621621621
$ws.Range("L18").Value = -11.450381679389
$ws.Range("M18").Value = -25.641025641025
$ws.Range("N18").Value = -89.004739336492

$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 50
$ws.Range("G19").Value = 38
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 520
$ws.Range("J19").Value = 677
$ws.Range("K19").Value = -23.190546528803
$ws.Range("L19").Value = -16.666666666666
$ws.Range("M19").Value = 60
$ws.Range("N19").Value = 25

$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 22
$ws.Range("G20").Value = 36
$ws.Range("H20").Value = -38.888888888888
$ws.Range("I20").Value = 256
$ws.Range("J20").Value = 325
$ws.Range("K20").Value = -21.230769230769
$ws.Range("L20").Value = -19.496855345911
$ws.Range("M20").Value = 116.949152542373
$ws.Range("N20").Value = -84.734645199761

$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = 4.761904761904
$ws.Range("F21").Value = 76
$ws.Range("G21").Value = 116
$ws.Range("H21").Value = -34.482758620689
$ws.Range("I21").Value = 1218
$ws.Range("J21").Value = 1527
$ws.Range("K21").Value = -20.235756385068
$ws.Range("L21").Value = -14.104372355430
$ws.Range("M21").Value = 38.882554161915
$ws.Range("N21").Value = -69.572820384711

$ws.Range("L22").Value = -41.666666666666

$ws.Range("D23").Value = 3
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = -87.5
$ws.Range("J23").Value = 65
$ws.Range("K23").Value = -24.615384615384
$ws.Range("M23").Value = 16.666666666666

$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = 41.176470588235
$ws.Range("F24").Value = 89
$ws.Range("G24").Value = 74
$ws.Range("H24").Value = 20.270270270270
$ws.Range("I24").Value = 1068
$ws.Range("J24").Value = 1014
$ws.Range("K24").Value = 5.325443786982
$ws.Range("L24").Value = 0.659754948162
$ws.Range("M24").Value = 42.780748663101

$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -71.428571428571
$ws.Range("F25").Value = 11
$ws.Range("G25").Value = 23
$ws.Range("H25").Value = -52.173913043478
$ws.Range("I25").Value = 287
$ws.Range("J25").Value = 371
$ws.Range("K25").Value = -22.641509433962
$ws.Range("L25").Value = -33.255813953488

$ws.Range("D26").Value = 7
$ws.Range("E26").Value = -14.285714285714
$ws.Range("F26").Value = 22
$ws.Range("G26").Value = 24
$ws.Range("H26").Value = -8.333333333333
$ws.Range("I26").Value = 345
$ws.Range("J26").Value = 336
$ws.Range("K26").Value = 2.678571428571
$ws.Range("L26").Value = 16.949152542372
$ws.Range("M26").Value = -8.244680851063

$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 300
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 32
$ws.Range("J27").Value = 22
$ws.Range("K27").Value = 45.454545454545
$ws.Range("L27").Value = 45.454545454545

$ws.Range("C28").Value = "0"
$ws.Range("G28").Value = "0"
$ws.Range("H28").Value = "***.*"

$ws.Range("C29").Value = 1
$ws.Range("F29").Value = 1
$ws.Range("I29").Value = 10
$ws.Range("K29").Value = 233.333333333333
$ws.Range("L29").Value = 66.666666666666
$ws.Range("M29").Value = 100
$ws.Range("N29").Value = -33.333333333333

$ws.Range("C30").Value = 1
$ws.Range("F30").Value = 1
$ws.Range("I30").Value = 7
$ws.Range("K30").Value = 133.333333333333
$ws.Range("L30").Value = 16.666666666666
$ws.Range("M30").Value = 40
$ws.Range("N30").Value = -41.666666666666

